$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.731.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "'1.854.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("D4").Value = "'1.012"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.74%  "
$ws.Range("D5").Value = "'319.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.47%  "
$ws.Range("D6").Value = "'1.010"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.84%  "
$ws.Range("D7").Value = "'0.4311"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.30%  "
$ws.Range("D8").Value = "'0.3752"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.69%  "
$ws.Range("D9").Value = "'0.07355"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.34%  "
$ws.Range("D10").Value = "'0.8787"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("D11").Value = "'21.65"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "'1.847.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("D13").Value = "'6.740"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("D14").Value = "'5.446"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("D15").Value = "'0.07138"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "'89.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.03%  "
$ws.Range("D17").Value = "'1.014"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.11%  "
$ws.Range("D18").Value = "'0.000009002"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("D19").Value = "'1.010"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.71%  "
$ws.Range("D20").Value = "'15.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.55%  "
$ws.Range("D21").Value = "'27.735.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("D22").Value = "'5.219"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.93%  "
$ws.Range("D23").Value = "'11.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.18%  "
$ws.Range("D24").Value = "'2.079.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.63%  "
$ws.Range("D25").Value = "'1.989"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.36%  "
$ws.Range("D26").Value = "'155.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.12%  "
$ws.Range("D27").Value = "'18.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("D28").Value = "'2.201"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.84%  "
$ws.Range("D29").Value = "'5.384"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").Value = "'119.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.90%  "
$ws.Range("D31").Value = "'0.08944"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.18%  "
$ws.Range("D32").Value = "'1.234"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("D33").Value = "'0.7794"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.96%  "
$ws.Range("D34").Value = "'4.564"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.69%  "
$ws.Range("D35").Value = "'2.925"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.47%  "
$ws.Range("D36").Value = "'1.011"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.76%  "
$ws.Range("D37").Value = "'1.131"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.47%  "
$ws.Range("D38").Value = "'0.05359"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("D39").Value = "'0.01979"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").Value = "'7.328"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.64%  "
$ws.Range("D41").Value = "'2.914"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.71%  "
$ws.Range("D42").Value = "'0.1695"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").Value = "'0.5144"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("D44").Value = "'8.839"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("D45").Value = "'10.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.78%  "
$ws.Range("D46").Value = "'108.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.81%  "
$ws.Range("D47").Value = "'0.4796"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("D48").Value = "'0.06482"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.85%  "
$ws.Range("D49").Value = "'1.693"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.76%  "
$ws.Range("D50").Value = "'1.011"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("D51").Value = "'1.850"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.16%  "
